{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the text replacements described by the diff: the date line and\n// 25 multiplication-fact cells inside the practice table.\nconst replacements = [\n  [\"2025-04-04 Friday\", \"2025-04-05 Saturday\"],\n  [\"519\u00d77=3633\", \"773\u00d76=4638\"],\n  [\"168\u00d78=1344\", \"417\u00d77=2919\"],\n  [\"256\u00d73=768\", \"688\u00d77=4816\"],\n  [\"545\u00d75=2725\", \"106\u00d79=954\"],\n  [\"760\u00d73=2280\", \"571\u00d79=5139\"],\n  [\"973\u00d74=3892\", \"213\u00d76=1278\"],\n  [\"631\u00d76=3786\", \"648\u00d72=1296\"],\n  [\"334\u00d77=2338\", \"561\u00d72=1122\"],\n  [\"286\u00d76=1716\", \"246\u00d73=738\"],\n  [\"320\u00d79=2880\", \"931\u00d72=1862\"],\n  [\"108\u00d72=216\", \"615\u00d72=1230\"],\n  [\"473\u00d79=4257\", \"628\u00d78=5024\"],\n  [\"581\u00d76=3486\", \"237\u00d79=2133\"],\n  [\"850\u00d78=6800\", \"864\u00d79=7776\"],\n  [\"450\u00d73=1350\", \"584\u00d77=4088\"],\n  [\"630\u00d74=2520\", \"410\u00d73=1230\"],\n  [\"138\u00d79=1242\", \"128\u00d73=384\"],\n  [\"268\u00d75=1340\", \"729\u00d74=2916\"],\n  [\"684\u00d74=2736\", \"521\u00d72=1042\"],\n  [\"406\u00d79=3654\", \"396\u00d79=3564\"],\n  [\"432\u00d72=864\", \"518\u00d76=3108\"],\n  [\"187\u00d73=561\", \"617\u00d72=1234\"],\n  [\"142\u00d75=710\", \"438\u00d75=2190\"],\n  [\"669\u00d76=4014\", \"718\u00d72=1436\"],\n  [\"365\u00d75=1825\", \"797\u00d77=5579\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Word COM interop script.\n# Applies the text replacements described by the diff: the date line and\n# 25 multiplication-fact cells inside the practice table.\n\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$replacements = @(\n    @(\"2025-04-04 Friday\", \"2025-04-05 Saturday\"),\n    @(\"519\u00d77=3633\", \"773\u00d76=4638\"),\n    @(\"168\u00d78=1344\", \"417\u00d77=2919\"),\n    @(\"256\u00d73=768\", \"688\u00d77=4816\"),\n    @(\"545\u00d75=2725\", \"106\u00d79=954\"),\n    @(\"760\u00d73=2280\", \"571\u00d79=5139\"),\n    @(\"973\u00d74=3892\", \"213\u00d76=1278\"),\n    @(\"631\u00d76=3786\", \"648\u00d72=1296\"),\n    @(\"334\u00d77=2338\", \"561\u00d72=1122\"),\n    @(\"286\u00d76=1716\", \"246\u00d73=738\"),\n    @(\"320\u00d79=2880\", \"931\u00d72=1862\"),\n    @(\"108\u00d72=216\", \"615\u00d72=1230\"),\n    @(\"473\u00d79=4257\", \"628\u00d78=5024\"),\n    @(\"581\u00d76=3486\", \"237\u00d79=2133\"),\n    @(\"850\u00d78=6800\", \"864\u00d79=7776\"),\n    @(\"450\u00d73=1350\", \"584\u00d77=4088\"),\n    @(\"630\u00d74=2520\", \"410\u00d73=1230\"),\n    @(\"138\u00d79=1242\", \"128\u00d73=384\"),\n    @(\"268\u00d75=1340\", \"729\u00d74=2916\"),\n    @(\"684\u00d74=2736\", \"521\u00d72=1042\"),\n    @(\"406\u00d79=3654\", \"396\u00d79=3564\"),\n    @(\"432\u00d72=864\", \"518\u00d76=3108\"),\n    @(\"187\u00d73=561\", \"617\u00d72=1234\"),\n    @(\"142\u00d75=710\", \"438\u00d75=2190\"),\n    @(\"669\u00d76=4014\", \"718\u00d72=1436\"),\n    @(\"365\u00d75=1825\", \"797\u00d77=5579\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll)\n}\n"}
